$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Mud"
$ws.Range("D4").Value = 1

$ws.Range("E7").Select()
